$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores free-form price text. Where the new price text
# would otherwise parse as a plain number (e.g. "600.22"), pre-format
# the cell as Text so Excel keeps it as a literal string, matching the
# thousands-dot-formatted text already used throughout this column.
$ws.Range("D2").Value = '68.680.19'
$ws.Range("E2").Value = '  +0.54%  '

$ws.Range("D3").Value = '2.710.18'
$ws.Range("E3").Value = '  +2.42%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.22'
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.43'
$ws.Range("E6").Value = '  +2.48%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").Value = '2.700.27'
$ws.Range("E9").Value = '  +2.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -4.19%  '

$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("E13").Value = '  +1.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.36'
$ws.Range("E14").Value = '  +1.05%  '

$ws.Range("D15").Value = '3.189.65'
$ws.Range("E15").Value = '  +2.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000187'
$ws.Range("E16").Value = '  -2.53%  '

$ws.Range("D17").Value = '68.738.09'
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").Value = '2.724.76'
$ws.Range("E18").Value = '  +3.03%  '

$ws.Range("E19").Value = '  +3.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  +3.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '367.85'
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.52'
$ws.Range("E22").Value = '  +2.94%  '

$ws.Range("E23").Value = '  +1.63%  '

$ws.Range("E24").Value = '  +2.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '75.37'
$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.13'
$ws.Range("E27").Value = '  +4.52%  '

$ws.Range("D28").Value = '2.806.65'
$ws.Range("E28").Value = '  +1.23%  '

$ws.Range("E29").Value = '  -1.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '586.19'
$ws.Range("E30").Value = '  +4.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.28'
$ws.Range("E32").Value = '  +3.26%  '

$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.95'
$ws.Range("E34").Value = '  +4.99%  '

$ws.Range("E35").Value = '  +5.94%  '

$ws.Range("E36").Value = '  +1.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.18'
$ws.Range("E38").Value = '  +4.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '161.16'
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  +2.25%  '

$ws.Range("E41").Value = '  +0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.64'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").Value = '0.0₆0319'
$ws.Range("E44").Value = '  -6.25%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '160.05'
$ws.Range("E46").Value = '  +0.46%  '

$ws.Range("E47").Value = '  +5.48%  '

$ws.Range("E48").Value = '  +6.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.614'
$ws.Range("E49").Value = '  +9.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.32'
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0789'
$ws.Range("E51").Value = '  +0.46%  '
